# JTown special collection picture addition first batch
#
# 1) Define the "Hyperlink" and "Unresolved Mention" character styles
#    (Word adds these automatically once a hyperlink is inserted /
#    collaboration features are touched).
# 2) Turn the freezingUrl line's URL into a real hyperlink - the run is
#    split so the hyperlink text also swallows the literal
#    "</freezingUrl" that precedes the closing angle bracket, matching
#    the source edit exactly.
# 3) Insert two blank paragraphs after that line.
# 4) Add a new paragraph containing the imgur picture link.

$d = $word.ActiveDocument

# --- 1) Styles -------------------------------------------------------
$hlStyle = $d.Styles.Add("Hyperlink", 2)
$hlStyle.BaseStyle = "DefaultParagraphFont"
$hlStyle.Priority = 99
$hlStyle.UnhideWhenUsed = $true
$hlStyle.Font.Underline = 1
$hlStyle.Font.TextColor.ObjectThemeColor = 10

$umStyle = $d.Styles.Add("Unresolved Mention", 2)
$umStyle.BaseStyle = "DefaultParagraphFont"
$umStyle.Priority = 99
$umStyle.UnhideWhenUsed = $true
$umStyle.Font.Color = 6053472

# --- 2) Split the freezingUrl paragraph and hyperlink the URL --------
$freezePara = $d.Paragraphs.Item(3)
$paraStart = $freezePara.Range.Start
$paraEnd = $freezePara.Range.End

$beforeText = "<freezingUrl>"
$linkText = "https://nchfp.uga.edu/how/freeze/apple.html</freezingUrl"
$afterText = ">"

$fullRange = $d.Range($paraStart, $paraEnd)
$fullRange.Text = $beforeText + $linkText + $afterText

$linkStart = $paraStart + $beforeText.Length
$linkEnd = $linkStart + $linkText.Length
$linkRange = $d.Range($linkStart, $linkEnd)

$d.Hyperlinks.Add($linkRange, "https://nchfp.uga.edu/how/freeze/apple.html", "", "", $linkText)

# --- 3) Two blank paragraphs after it ---------------------------------
$freezePara = $d.Paragraphs.Item(3)
$freezePara.Range.InsertParagraphAfter()

$blankPara2 = $d.Paragraphs.Item(4)
$blankPara2.Range.InsertParagraphAfter()

# --- 4) New paragraph with the imgur picture link ---------------------
$picPara = $d.Paragraphs.Item(5)
$picPara.Range.InsertParagraphAfter()

$imgPara = $d.Paragraphs.Item(6)
$imgPara.Range.Text = "https://i.imgur.com/ke4SM3b.jpg"
